# Update column G ("K" = strikeouts) values for rows 2-14 on the active
# sheet, per the regenerated save_data (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 2
    6  = 3
    7  = 1
    8  = 2
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
